$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" field: 9/20/2017 -> 8/18/2018
#    Appears once on the slide master and once on every slide layout
#    (the "Date Placeholder" shape).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "9/20/2017") {
                $shp.TextFrame.TextRange.Text = "8/18/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 21 ("Техника инстанцирования"), third shape: the two
#    Consolas code blocks had their signature split across two runs
#    ("max<int> " + "(int x, int y) {", and "max<double> " + "(double ")
#    -- re-typing/merging them into a single run each (formatting is
#    identical on both sides of the former split, so this just
#    coalesces runs, it doesn't change the visible text).
# ---------------------------------------------------------------------
function Merge-AdjacentRuns($textRange, [string]$firstPart, [string]$secondPart) {
    $full = $textRange.Text
    $needle = $firstPart + $secondPart
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0) {
        $span = $textRange.Characters($idx + 1, $needle.Length)
        # Re-assigning the identical text forces the run that spans the
        # old run boundary to be rebuilt as a single run (same rPr on
        # both sides, so formatting is unaffected).
        $span.Text = $span.Text
    }
}

$slide21 = $p.Slides.Item(21)
$codeShape = $slide21.Shapes.Item(3)
$tr = $codeShape.TextFrame.TextRange

Merge-AdjacentRuns $tr "max<int> " "(int x, int y) {"
Merge-AdjacentRuns $tr "max<double> " "(double "
